$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# 1. PBIReports sheet: rename the "PBIReportId" column to
#    "PBIReportIdFieldName" and replace each row's GUID value with the
#    new constant field name "Finance Report Id".
# ---------------------------------------------------------------------
$wsReports = $wb.Worksheets.Item("PBIReports")
$wsReports.Range("G1").Value = "PBIReportIdFieldName"
for ($r = 2; $r -le 7; $r++) {
    $wsReports.Cells.Item($r, 7).Value = "Finance Report Id"
}

# ---------------------------------------------------------------------
# 2. Namespace sheet: rename the namespace value.
# ---------------------------------------------------------------------
$wsNamespace = $wb.Worksheets.Item("Namespace")
$wsNamespace.Range("A2").Value = "Microsoft.Finance.PowerBIReports"

# ---------------------------------------------------------------------
# 3. View/selection updates to match the report-selection feature.
# ---------------------------------------------------------------------

# PBIReports: move the active selection to H1 (no tab activation here).
$wsReports.Select()
$wsReports.Range("H1").Select()

# Namespace: becomes the active tab, with A3 selected.
$wsNamespace.Activate()
$wsNamespace.Range("A3").Select()

Write-Output "done"
